# Apply crypto price/volume updates from the Fri Oct 20 19:29:52 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.440.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "'1.607.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'212.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("E6").Value = "  +7.05%  "
$ws.Range("D8").Value = "'26.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").Value = "'43.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "'0.0910"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").Value = "'1.838.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "'1.603.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "'29.482.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "'63.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").Value = "'240.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.46%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.82%  "
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'3.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("D24").Value = "'9.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'154.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("D28").Value = "'15.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "'3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("D34").Value = "'3.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("D35").Value = "'1.414.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("D38").Value = "'2.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.11%  "
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +6.40%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.797"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("D46").Value = "'52.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +22.02%  "
$ws.Range("D47").Value = "'65.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "'1.747.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "'86.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.88%  "
